$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.676.71"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.533.94"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'288.25"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.3934"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("D8").Value = "'0.3152"
$ws.Range("E8").Value = "  -3.02%  "
$ws.Range("D9").Value = "'42.41"
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").Value = "'0.07154"
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("D11").Value = "'1.045"
$ws.Range("E11").Value = "  -7.41%  "
$ws.Range("D12").Value = "'0.9997"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "'5.644"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").Value = "'18.53"
$ws.Range("E14").Value = "  -4.82%  "
$ws.Range("D15").Value = "'6.588"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("D16").Value = "1.535.93"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "'0.00001090"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "'0.06593"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "'83.27"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "'6.112"
$ws.Range("E21").Value = "  -4.66%  "
$ws.Range("D22").Value = "'15.40"
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("D23").Value = "'10.80"
$ws.Range("E23").Value = "  -6.07%  "
$ws.Range("D24").Value = "'2.357"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "21.678.79"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").Value = "'2.345"
$ws.Range("E26").Value = "  -7.96%  "
$ws.Range("D27").Value = "'147.99"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").Value = "'18.32"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("D29").Value = "'4.836"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").Value = "1.709.05"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "'116.98"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").Value = "'5.879"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").Value = "'0.9423"
$ws.Range("E33").Value = "  -15.68%  "
$ws.Range("D34").Value = "'0.08129"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "'8.507"
$ws.Range("E35").Value = "  -8.47%  "
$ws.Range("D36").Value = "'5.117"
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("D37").Value = "'0.05996"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("D38").Value = "'0.02201"
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("E39").Value = "  -14.55%  "
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.171"
$ws.Range("E41").Value = "  -4.29%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'10.94"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "'0.9998"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "'0.5748"
$ws.Range("E44").Value = "  -3.54%  "
$ws.Range("D45").Value = "'13.02"
$ws.Range("E45").Value = "  -4.09%  "
$ws.Range("D46").Value = "'3.704"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "'0.5481"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "'116.37"
$ws.Range("E49").Value = "  -2.70%  "
$ws.Range("D50").Value = "'1.864"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("D51").Value = "'0.06691"
$ws.Range("E51").Value = "  -3.00%  "
